$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(315).Delete()
